$d = $word.ActiveDocument

# Locate the "Angles capturing:" run so we can split its paragraph right
# after the text and before the existing _GoBack bookmark.
$rng = $d.Content
$found = $rng.Find.Execute("Angles capturing:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the point right after the found text (still before the
    # bookmark that lives at the end of the same paragraph).
    $rng.Collapse(0)

    # Insert two paragraph breaks at that point: this turns the single
    # paragraph "Angles capturing:<bookmark>" into three paragraphs -
    # "Angles capturing:", an empty one, and the one holding the bookmark -
    # all sharing the same (jc=both) paragraph formatting.
    $rng.InsertBefore("`r")
    $rng.InsertBefore("`r")
}
